$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, date range) ---
$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Crime data table updates (rows 14-29) ---
$ws.Range("C14").Value = "'0"
$ws.Range("A36").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("N14").Value = -82.857142857142
$ws.Range("C36").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("K36").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = "'0"
$ws.Range("A36").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = -15
$ws.Range("N15").Value = -67.924528301886
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -5.882352941176
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 131
$ws.Range("K16").Value = -6.870229007633
$ws.Range("L16").Value = 46.987951807228
$ws.Range("M16").Value = -29.885057471264
$ws.Range("N16").Value = -76.893939393939
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -25
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 15.384615384615
$ws.Range("I17").Value = 278
$ws.Range("J17").Value = 326
$ws.Range("K17").Value = -14.723926380368
$ws.Range("L17").Value = -6.711409395973
$ws.Range("M17").Value = 47.089947089947
$ws.Range("N17").Value = -54.125412541254
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 22
$ws.Range("H18").Value = 144.444444444444
$ws.Range("I18").Value = 112
$ws.Range("J18").Value = 124
$ws.Range("K18").Value = -9.677419354838
$ws.Range("L18").Value = 13.131313131313
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = -80.521739130434
$ws.Range("C19").Value = 12
$ws.Range("E19").Value = 140
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 125
$ws.Range("I19").Value = 250
$ws.Range("J19").Value = 206
$ws.Range("K19").Value = 21.359223300970
$ws.Range("L19").Value = 33.689839572192
$ws.Range("M19").Value = 37.362637362637
$ws.Range("N19").Value = 9.649122807017
$ws.Range("I20").Value = 65
$ws.Range("K20").Value = 30
$ws.Range("L20").Value = 124.137931034483
$ws.Range("M20").Value = 124.137931034483
$ws.Range("N20").Value = -59.627329192546
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 140
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = 50.537634408602
$ws.Range("I21").Value = 850
$ws.Range("J21").Value = 860
$ws.Range("K21").Value = -1.162790697674
$ws.Range("L21").Value = 18.715083798882
$ws.Range("M21").Value = 27.627627627627
$ws.Range("N21").Value = -61.116193961573
$ws.Range("F22").Value = "'0"
$ws.Range("A36").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = 100
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 25
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 27.777777777777
$ws.Range("I23").Value = 141
$ws.Range("J23").Value = 134
$ws.Range("K23").Value = 5.223880597014
$ws.Range("L23").Value = 14.634146341463
$ws.Range("M23").Value = 46.875
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 175
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = 65.384615384615
$ws.Range("I24").Value = 519
$ws.Range("J24").Value = 528
$ws.Range("K24").Value = -1.704545454545
$ws.Range("L24").Value = 15.077605321507
$ws.Range("M24").Value = 39.516129032258
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -10
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 10
$ws.Range("I25").Value = 362
$ws.Range("J25").Value = 350
$ws.Range("K25").Value = 3.428571428571
$ws.Range("L25").Value = 4.322766570605
$ws.Range("M25").Value = -32.962962962963
$ws.Range("C36").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1
$ws.Range("C36").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("K36").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 29
$ws.Range("K26").Value = 3.448275862068
$ws.Range("L26").Value = 30.434782608695
$ws.Range("C27").Value = 1
$ws.Range("C36").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("K36").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = -28.846153846153
$ws.Range("L27").Value = -22.916666666666
$ws.Range("C28").Value = "'0"
$ws.Range("A36").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("F28").Value = 4
$ws.Range("L28").Value = -41.666666666666
$ws.Range("M28").Value = -8.695652173913
$ws.Range("N28").Value = -72.727272727272
$ws.Range("C29").Value = "'0"
$ws.Range("A36").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 3
$ws.Range("L29").Value = -41.935483870967
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -74.647887323943
